# Insert a new row at row 21. This shifts the existing rows 21-179 down to 22-180,
# preserving all of their data, and grows the used range from A1:R179 to A1:R180.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with the same static data as the row that was
# pushed down to row 22 (i.e. what used to be row 21), but with the new Fecha
# (column D) and Volumen (column J) values from the commit.
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = "2021-10-04"
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100114014
$ws.Range("G21").Value = "Betarraga"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 650
$ws.Range("L21").Value = 650
$ws.Range("M21").Value = 650
$ws.Range("N21").Value = "`$/paquete 5 unidades"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 130
$ws.Range("Q21").Value = 5
$ws.Range("R21").Value = "Hortaliza"
